$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("PIR")
$ws.Range("A14:A28").NumberFormat = "@"

$ws.Cells.Item(14, 1).Value = "2026-02-01"
$ws.Cells.Item(14, 2).Value = "19:53:40"
$ws.Cells.Item(14, 3).Value = "19:00"
$ws.Cells.Item(14, 4).Value = "Bathroom"
$ws.Cells.Item(14, 5).Value = "Motion Detected"
$ws.Cells.Item(14, 6).Value = "Active"

$ws.Cells.Item(15, 1).Value = "2026-02-01"
$ws.Cells.Item(15, 2).Value = "19:53:41"
$ws.Cells.Item(15, 3).Value = "19:00"
$ws.Cells.Item(15, 4).Value = "Bathroom"
$ws.Cells.Item(15, 5).Value = "No Motion"
$ws.Cells.Item(15, 6).Value = "Inactive"

$ws.Cells.Item(16, 1).Value = "2026-02-01"
$ws.Cells.Item(16, 2).Value = "19:53:41"
$ws.Cells.Item(16, 3).Value = "19:00"
$ws.Cells.Item(16, 4).Value = "Bathroom"
$ws.Cells.Item(16, 5).Value = "No Motion"
$ws.Cells.Item(16, 6).Value = "Inactive"

$ws.Cells.Item(17, 1).Value = "2026-02-01"
$ws.Cells.Item(17, 2).Value = "19:53:46"
$ws.Cells.Item(17, 3).Value = "19:00"
$ws.Cells.Item(17, 4).Value = "Bathroom"
$ws.Cells.Item(17, 5).Value = "No Motion"
$ws.Cells.Item(17, 6).Value = "Inactive"

$ws.Cells.Item(18, 1).Value = "2026-02-01"
$ws.Cells.Item(18, 2).Value = "19:53:51"
$ws.Cells.Item(18, 3).Value = "19:00"
$ws.Cells.Item(18, 4).Value = "Bathroom"
$ws.Cells.Item(18, 5).Value = "No Motion"
$ws.Cells.Item(18, 6).Value = "Inactive"

$ws.Cells.Item(19, 1).Value = "2026-02-01"
$ws.Cells.Item(19, 2).Value = "19:53:56"
$ws.Cells.Item(19, 3).Value = "19:00"
$ws.Cells.Item(19, 4).Value = "Bathroom"
$ws.Cells.Item(19, 5).Value = "No Motion"
$ws.Cells.Item(19, 6).Value = "Inactive"

$ws.Cells.Item(20, 1).Value = "2026-02-01"
$ws.Cells.Item(20, 2).Value = "19:54:01"
$ws.Cells.Item(20, 3).Value = "19:00"
$ws.Cells.Item(20, 4).Value = "Bathroom"
$ws.Cells.Item(20, 5).Value = "No Motion"
$ws.Cells.Item(20, 6).Value = "Inactive"

$ws.Cells.Item(21, 1).Value = "2026-02-01"
$ws.Cells.Item(21, 2).Value = "19:54:06"
$ws.Cells.Item(21, 3).Value = "19:00"
$ws.Cells.Item(21, 4).Value = "Bathroom"
$ws.Cells.Item(21, 5).Value = "No Motion"
$ws.Cells.Item(21, 6).Value = "Inactive"

$ws.Cells.Item(22, 1).Value = "2026-02-01"
$ws.Cells.Item(22, 2).Value = "19:54:11"
$ws.Cells.Item(22, 3).Value = "19:00"
$ws.Cells.Item(22, 4).Value = "Bathroom"
$ws.Cells.Item(22, 5).Value = "No Motion"
$ws.Cells.Item(22, 6).Value = "Inactive"

$ws.Cells.Item(23, 1).Value = "2026-02-01"
$ws.Cells.Item(23, 2).Value = "19:54:16"
$ws.Cells.Item(23, 3).Value = "19:00"
$ws.Cells.Item(23, 4).Value = "Bathroom"
$ws.Cells.Item(23, 5).Value = "No Motion"
$ws.Cells.Item(23, 6).Value = "Inactive"

$ws.Cells.Item(24, 1).Value = "2026-02-01"
$ws.Cells.Item(24, 2).Value = "19:54:21"
$ws.Cells.Item(24, 3).Value = "19:00"
$ws.Cells.Item(24, 4).Value = "Bathroom"
$ws.Cells.Item(24, 5).Value = "No Motion"
$ws.Cells.Item(24, 6).Value = "Inactive"

$ws.Cells.Item(25, 1).Value = "2026-02-01"
$ws.Cells.Item(25, 2).Value = "19:54:26"
$ws.Cells.Item(25, 3).Value = "19:00"
$ws.Cells.Item(25, 4).Value = "Bathroom"
$ws.Cells.Item(25, 5).Value = "No Motion"
$ws.Cells.Item(25, 6).Value = "Inactive"

$ws.Cells.Item(26, 1).Value = "2026-02-01"
$ws.Cells.Item(26, 2).Value = "19:54:31"
$ws.Cells.Item(26, 3).Value = "19:00"
$ws.Cells.Item(26, 4).Value = "Bathroom"
$ws.Cells.Item(26, 5).Value = "No Motion"
$ws.Cells.Item(26, 6).Value = "Inactive"

$ws.Cells.Item(27, 1).Value = "2026-02-01"
$ws.Cells.Item(27, 2).Value = "19:54:37"
$ws.Cells.Item(27, 3).Value = "19:00"
$ws.Cells.Item(27, 4).Value = "Bathroom"
$ws.Cells.Item(27, 5).Value = "No Motion"
$ws.Cells.Item(27, 6).Value = "Inactive"

$ws.Cells.Item(28, 1).Value = "2026-02-01"
$ws.Cells.Item(28, 2).Value = "19:54:37"
$ws.Cells.Item(28, 3).Value = "19:00"
$ws.Cells.Item(28, 4).Value = "Bathroom"
$ws.Cells.Item(28, 5).Value = "Motion Detected"
$ws.Cells.Item(28, 6).Value = "Active"

$ws = $wb.Worksheets.Item("Humidity")
$ws.Range("A11:A23").NumberFormat = "@"
$ws.Range("E11:E23").NumberFormat = "@"

$ws.Cells.Item(11, 1).Value = "2026-02-01"
$ws.Cells.Item(11, 2).Value = "19:53:40"
$ws.Cells.Item(11, 3).Value = "19:00"
$ws.Cells.Item(11, 4).Value = "Bathroom"
$ws.Cells.Item(11, 5).Value = "78.6%"
$ws.Cells.Item(11, 6).Value = "Active"

$ws.Cells.Item(12, 1).Value = "2026-02-01"
$ws.Cells.Item(12, 2).Value = "19:53:41"
$ws.Cells.Item(12, 3).Value = "19:00"
$ws.Cells.Item(12, 4).Value = "Bathroom"
$ws.Cells.Item(12, 5).Value = "78.4%"
$ws.Cells.Item(12, 6).Value = "Active"

$ws.Cells.Item(13, 1).Value = "2026-02-01"
$ws.Cells.Item(13, 2).Value = "19:53:45"
$ws.Cells.Item(13, 3).Value = "19:00"
$ws.Cells.Item(13, 4).Value = "Bathroom"
$ws.Cells.Item(13, 5).Value = "77.8%"
$ws.Cells.Item(13, 6).Value = "Active"

$ws.Cells.Item(14, 1).Value = "2026-02-01"
$ws.Cells.Item(14, 2).Value = "19:53:50"
$ws.Cells.Item(14, 3).Value = "19:00"
$ws.Cells.Item(14, 4).Value = "Bathroom"
$ws.Cells.Item(14, 5).Value = "79.0%"
$ws.Cells.Item(14, 6).Value = "Active"

$ws.Cells.Item(15, 1).Value = "2026-02-01"
$ws.Cells.Item(15, 2).Value = "19:53:55"
$ws.Cells.Item(15, 3).Value = "19:00"
$ws.Cells.Item(15, 4).Value = "Bathroom"
$ws.Cells.Item(15, 5).Value = "77.3%"
$ws.Cells.Item(15, 6).Value = "Active"

$ws.Cells.Item(16, 1).Value = "2026-02-01"
$ws.Cells.Item(16, 2).Value = "19:54:00"
$ws.Cells.Item(16, 3).Value = "19:00"
$ws.Cells.Item(16, 4).Value = "Bathroom"
$ws.Cells.Item(16, 5).Value = "77.9%"
$ws.Cells.Item(16, 6).Value = "Active"

$ws.Cells.Item(17, 1).Value = "2026-02-01"
$ws.Cells.Item(17, 2).Value = "19:54:05"
$ws.Cells.Item(17, 3).Value = "19:00"
$ws.Cells.Item(17, 4).Value = "Bathroom"
$ws.Cells.Item(17, 5).Value = "77.1%"
$ws.Cells.Item(17, 6).Value = "Active"

$ws.Cells.Item(18, 1).Value = "2026-02-01"
$ws.Cells.Item(18, 2).Value = "19:54:10"
$ws.Cells.Item(18, 3).Value = "19:00"
$ws.Cells.Item(18, 4).Value = "Bathroom"
$ws.Cells.Item(18, 5).Value = "77.8%"
$ws.Cells.Item(18, 6).Value = "Active"

$ws.Cells.Item(19, 1).Value = "2026-02-01"
$ws.Cells.Item(19, 2).Value = "19:54:15"
$ws.Cells.Item(19, 3).Value = "19:00"
$ws.Cells.Item(19, 4).Value = "Bathroom"
$ws.Cells.Item(19, 5).Value = "76.9%"
$ws.Cells.Item(19, 6).Value = "Active"

$ws.Cells.Item(20, 1).Value = "2026-02-01"
$ws.Cells.Item(20, 2).Value = "19:54:20"
$ws.Cells.Item(20, 3).Value = "19:00"
$ws.Cells.Item(20, 4).Value = "Bathroom"
$ws.Cells.Item(20, 5).Value = "78.1%"
$ws.Cells.Item(20, 6).Value = "Active"

$ws.Cells.Item(21, 1).Value = "2026-02-01"
$ws.Cells.Item(21, 2).Value = "19:54:25"
$ws.Cells.Item(21, 3).Value = "19:00"
$ws.Cells.Item(21, 4).Value = "Bathroom"
$ws.Cells.Item(21, 5).Value = "77.3%"
$ws.Cells.Item(21, 6).Value = "Active"

$ws.Cells.Item(22, 1).Value = "2026-02-01"
$ws.Cells.Item(22, 2).Value = "19:54:31"
$ws.Cells.Item(22, 3).Value = "19:00"
$ws.Cells.Item(22, 4).Value = "Bathroom"
$ws.Cells.Item(22, 5).Value = "78.5%"
$ws.Cells.Item(22, 6).Value = "Active"

$ws.Cells.Item(23, 1).Value = "2026-02-01"
$ws.Cells.Item(23, 2).Value = "19:54:36"
$ws.Cells.Item(23, 3).Value = "19:00"
$ws.Cells.Item(23, 4).Value = "Bathroom"
$ws.Cells.Item(23, 5).Value = "77.2%"
$ws.Cells.Item(23, 6).Value = "Active"

$ws = $wb.Worksheets.Item("Temperature")
$ws.Range("A11:A23").NumberFormat = "@"

$ws.Cells.Item(11, 1).Value = "2026-02-01"
$ws.Cells.Item(11, 2).Value = "19:53:40"
$ws.Cells.Item(11, 3).Value = "19:00"
$ws.Cells.Item(11, 4).Value = "Bathroom"
$ws.Cells.Item(11, 5).Value = "25.3C"
$ws.Cells.Item(11, 6).Value = "Active"

$ws.Cells.Item(12, 1).Value = "2026-02-01"
$ws.Cells.Item(12, 2).Value = "19:53:41"
$ws.Cells.Item(12, 3).Value = "19:00"
$ws.Cells.Item(12, 4).Value = "Bathroom"
$ws.Cells.Item(12, 5).Value = "25.3C"
$ws.Cells.Item(12, 6).Value = "Active"

$ws.Cells.Item(13, 1).Value = "2026-02-01"
$ws.Cells.Item(13, 2).Value = "19:53:45"
$ws.Cells.Item(13, 3).Value = "19:00"
$ws.Cells.Item(13, 4).Value = "Bathroom"
$ws.Cells.Item(13, 5).Value = "25.3C"
$ws.Cells.Item(13, 6).Value = "Active"

$ws.Cells.Item(14, 1).Value = "2026-02-01"
$ws.Cells.Item(14, 2).Value = "19:53:50"
$ws.Cells.Item(14, 3).Value = "19:00"
$ws.Cells.Item(14, 4).Value = "Bathroom"
$ws.Cells.Item(14, 5).Value = "25.3C"
$ws.Cells.Item(14, 6).Value = "Active"

$ws.Cells.Item(15, 1).Value = "2026-02-01"
$ws.Cells.Item(15, 2).Value = "19:53:55"
$ws.Cells.Item(15, 3).Value = "19:00"
$ws.Cells.Item(15, 4).Value = "Bathroom"
$ws.Cells.Item(15, 5).Value = "25.3C"
$ws.Cells.Item(15, 6).Value = "Active"

$ws.Cells.Item(16, 1).Value = "2026-02-01"
$ws.Cells.Item(16, 2).Value = "19:54:00"
$ws.Cells.Item(16, 3).Value = "19:00"
$ws.Cells.Item(16, 4).Value = "Bathroom"
$ws.Cells.Item(16, 5).Value = "25.3C"
$ws.Cells.Item(16, 6).Value = "Active"

$ws.Cells.Item(17, 1).Value = "2026-02-01"
$ws.Cells.Item(17, 2).Value = "19:54:05"
$ws.Cells.Item(17, 3).Value = "19:00"
$ws.Cells.Item(17, 4).Value = "Bathroom"
$ws.Cells.Item(17, 5).Value = "25.3C"
$ws.Cells.Item(17, 6).Value = "Active"

$ws.Cells.Item(18, 1).Value = "2026-02-01"
$ws.Cells.Item(18, 2).Value = "19:54:10"
$ws.Cells.Item(18, 3).Value = "19:00"
$ws.Cells.Item(18, 4).Value = "Bathroom"
$ws.Cells.Item(18, 5).Value = "25.3C"
$ws.Cells.Item(18, 6).Value = "Active"

$ws.Cells.Item(19, 1).Value = "2026-02-01"
$ws.Cells.Item(19, 2).Value = "19:54:16"
$ws.Cells.Item(19, 3).Value = "19:00"
$ws.Cells.Item(19, 4).Value = "Bathroom"
$ws.Cells.Item(19, 5).Value = "25.3C"
$ws.Cells.Item(19, 6).Value = "Active"

$ws.Cells.Item(20, 1).Value = "2026-02-01"
$ws.Cells.Item(20, 2).Value = "19:54:20"
$ws.Cells.Item(20, 3).Value = "19:00"
$ws.Cells.Item(20, 4).Value = "Bathroom"
$ws.Cells.Item(20, 5).Value = "25.2C"
$ws.Cells.Item(20, 6).Value = "Active"

$ws.Cells.Item(21, 1).Value = "2026-02-01"
$ws.Cells.Item(21, 2).Value = "19:54:25"
$ws.Cells.Item(21, 3).Value = "19:00"
$ws.Cells.Item(21, 4).Value = "Bathroom"
$ws.Cells.Item(21, 5).Value = "25.3C"
$ws.Cells.Item(21, 6).Value = "Active"

$ws.Cells.Item(22, 1).Value = "2026-02-01"
$ws.Cells.Item(22, 2).Value = "19:54:31"
$ws.Cells.Item(22, 3).Value = "19:00"
$ws.Cells.Item(22, 4).Value = "Bathroom"
$ws.Cells.Item(22, 5).Value = "25.3C"
$ws.Cells.Item(22, 6).Value = "Active"

$ws.Cells.Item(23, 1).Value = "2026-02-01"
$ws.Cells.Item(23, 2).Value = "19:54:37"
$ws.Cells.Item(23, 3).Value = "19:00"
$ws.Cells.Item(23, 4).Value = "Bathroom"
$ws.Cells.Item(23, 5).Value = "25.3C"
$ws.Cells.Item(23, 6).Value = "Active"
